$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.376.71'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.105.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.98%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.101.20'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.446'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.26'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('E12').Value = '  +2.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.640.52'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('E14').Value = '  +3.15%  '
$ws.Range('E15').Value = '  -2.40%  '
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '57.471.33'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.102.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.43'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '346.48'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.80'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '67.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.64%  '
$ws.Range('E26').Value = '  -1.88%  '
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0894'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.42%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.37'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.15%  '
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.03'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.86'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.93'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.29%  '
$ws.Range('E36').Value = '  -2.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.97'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.06'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.80'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.58%  '
$ws.Range('E40').Value = '  -1.20%  '
$ws.Range('E41').Value = '  +5.99%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.61'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.64%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0660'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.699'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.144.25'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '36.58'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.354.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.45%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('E49').Value = '  +3.27%  '
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('E51').Value = '  -0.37%  '
